$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Median Value (column C) and Tier (column D) for each school row.
# Values were recalculated relative to the median AFTER merging with zip/census tract data.
$ws.Range("C2").Value = 1.049645390070922
$ws.Range("D2").Value = "4th Tier"
$ws.Range("C3").Value = 0.9148936170212766
$ws.Range("D3").Value = "Below Median"
$ws.Range("C4").Value = 1.521276595744681
$ws.Range("D4").Value = "1st Tier"
$ws.Range("C5").Value = 1.035460992907801
$ws.Range("D5").Value = "4th Tier"
$ws.Range("C6").Value = 1.521276595744681
$ws.Range("D6").Value = "1st Tier"
$ws.Range("C7").Value = 1.893617021276596
$ws.Range("D7").Value = "1st Tier"
$ws.Range("C8").Value = 1.148936170212766
$ws.Range("D8").Value = "3rd Tier"
$ws.Range("C9").Value = 1.893617021276596
$ws.Range("D9").Value = "1st Tier"
$ws.Range("C10").Value = 0.1170212765957447
$ws.Range("D10").Value = "Below Median"
$ws.Range("C11").Value = 0.1170212765957447
$ws.Range("D11").Value = "Below Median"
$ws.Range("C12").Value = 0.1170212765957447
$ws.Range("D12").Value = "Below Median"
$ws.Range("C13").Value = 0.1170212765957447
$ws.Range("D13").Value = "Below Median"
$ws.Range("C14").Value = 0.1702127659574468
$ws.Range("D14").Value = "Below Median"
$ws.Range("C15").Value = 0.1702127659574468
$ws.Range("D15").Value = "Below Median"
$ws.Range("C16").Value = 0.425531914893617
$ws.Range("D16").Value = "Below Median"
$ws.Range("C17").Value = 0.5212765957446809
$ws.Range("D17").Value = "Below Median"
$ws.Range("C18").Value = 0.6382978723404256
$ws.Range("D18").Value = "Below Median"
$ws.Range("C19").Value = 0.4680851063829787
$ws.Range("D19").Value = "Below Median"
$ws.Range("C20").Value = 0.4680851063829787
$ws.Range("D20").Value = "Below Median"
$ws.Range("C21").Value = 0.5531914893617021
$ws.Range("D21").Value = "Below Median"
$ws.Range("C22").Value = 0.2553191489361702
$ws.Range("D22").Value = "Below Median"
$ws.Range("C23").Value = 0.2553191489361702
$ws.Range("D23").Value = "Below Median"
$ws.Range("C24").Value = 1.340425531914894
$ws.Range("D24").Value = "2nd Tier"
$ws.Range("C25").Value = 1.453900709219858
$ws.Range("D25").Value = "2nd Tier"
$ws.Range("C26").Value = 1.606382978723404
$ws.Range("D26").Value = "1st Tier"
$ws.Range("C27").Value = 1.606382978723404
$ws.Range("D27").Value = "1st Tier"
$ws.Range("C28").Value = 1.627659574468085
$ws.Range("D28").Value = "1st Tier"
$ws.Range("C29").Value = 1.134751773049645
$ws.Range("D29").Value = "3rd Tier"
$ws.Range("C30").Value = 1.390070921985816
$ws.Range("D30").Value = "2nd Tier"
$ws.Range("C31").Value = 0.925531914893617
$ws.Range("D31").Value = "Below Median"
$ws.Range("C32").Value = 1.340425531914894
$ws.Range("D32").Value = "2nd Tier"
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = "4th Tier"
$ws.Range("C34").Value = 0.8404255319148937
$ws.Range("D34").Value = "Below Median"
$ws.Range("C35").Value = 0.7943262411347518
$ws.Range("D35").Value = "Below Median"
$ws.Range("C36").Value = 1.326241134751773
$ws.Range("D36").Value = "3rd Tier"
$ws.Range("C37").Value = 1.361702127659574
$ws.Range("D37").Value = "2nd Tier"
$ws.Range("C38").Value = 1.361702127659574
$ws.Range("D38").Value = "2nd Tier"
$ws.Range("C39").Value = 1.319148936170213
$ws.Range("D39").Value = "3rd Tier"
$ws.Range("C40").Value = 1.326241134751773
$ws.Range("D40").Value = "2nd Tier"
$ws.Range("C41").Value = 1.361702127659574
$ws.Range("D41").Value = "2nd Tier"
$ws.Range("C42").Value = 1.014184397163121
$ws.Range("D42").Value = "4th Tier"
$ws.Range("C43").Value = 1.340425531914894
$ws.Range("D43").Value = "2nd Tier"
$ws.Range("C44").Value = 1.056737588652482
$ws.Range("D44").Value = "4th Tier"
$ws.Range("C45").Value = 1.049645390070922
$ws.Range("D45").Value = "4th Tier"
$ws.Range("C46").Value = 0.3829787234042553
$ws.Range("D46").Value = "Below Median"
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = "4th Tier"
$ws.Range("C48").Value = 0.8457446808510638
$ws.Range("D48").Value = "Below Median"
$ws.Range("C49").Value = 0.851063829787234
$ws.Range("D49").Value = "Below Median"
$ws.Range("C50").Value = 0.7943262411347518
$ws.Range("D50").Value = "Below Median"
$ws.Range("C51").Value = 1.049645390070922
$ws.Range("D51").Value = "4th Tier"
$ws.Range("C52").Value = 1
$ws.Range("D52").Value = "4th Tier"
$ws.Range("C53").Value = 1.308510638297872
$ws.Range("D53").Value = "3rd Tier"
$ws.Range("C54").Value = 0.4468085106382979
$ws.Range("D54").Value = "Below Median"
$ws.Range("C55").Value = 1.095744680851064
$ws.Range("D55").Value = "3rd Tier"
$ws.Range("C56").Value = 1.106382978723404
$ws.Range("D56").Value = "3rd Tier"
$ws.Range("C57").Value = 0.4609929078014184
$ws.Range("D57").Value = "Below Median"
$ws.Range("C58").Value = 0.3333333333333333
$ws.Range("D58").Value = "Below Median"
$ws.Range("C59").Value = 0.9787234042553191
$ws.Range("D59").Value = "Below Median"
$ws.Range("C60").Value = 1.24468085106383
$ws.Range("D60").Value = "3rd Tier"
$ws.Range("C61").Value = 1.148936170212766
$ws.Range("D61").Value = "3rd Tier"
$ws.Range("C62").Value = 1.319148936170213
$ws.Range("D62").Value = "3rd Tier"
$ws.Range("C63").Value = 1.014184397163121
$ws.Range("D63").Value = "4th Tier"
$ws.Range("C64").Value = 0.3333333333333333
$ws.Range("D64").Value = "Below Median"
$ws.Range("C65").Value = 1
$ws.Range("D65").Value = "4th Tier"
$ws.Range("C66").Value = 0.6879432624113476
$ws.Range("D66").Value = "Below Median"
$ws.Range("C67").Value = 0.4609929078014184
$ws.Range("D67").Value = "Below Median"
$ws.Range("C68").Value = 1.453900709219858
$ws.Range("D68").Value = "2nd Tier"
$ws.Range("C69").Value = 0.148936170212766
$ws.Range("D69").Value = "Below Median"
$ws.Range("C70").Value = 0.6808510638297872
$ws.Range("D70").Value = "Below Median"
$ws.Range("C71").Value = 1.453900709219858
$ws.Range("D71").Value = "2nd Tier"
$ws.Range("C72").Value = 0.1170212765957447
$ws.Range("D72").Value = "Below Median"
$ws.Range("C73").Value = 0.3191489361702128
$ws.Range("D73").Value = "Below Median"
$ws.Range("C74").Value = 0.8191489361702128
$ws.Range("D74").Value = "Below Median"
$ws.Range("C75").Value = 1.521276595744681
$ws.Range("D75").Value = "1st Tier"
$ws.Range("C76").Value = 1.453900709219858
$ws.Range("D76").Value = "1st Tier"
$ws.Range("C77").Value = 1.085106382978723
$ws.Range("D77").Value = "3rd Tier"
$ws.Range("C78").Value = 1.524822695035461
$ws.Range("D78").Value = "1st Tier"
$ws.Range("C79").Value = 0.9148936170212766
$ws.Range("D79").Value = "Below Median"
$ws.Range("C80").Value = 0.6666666666666666
$ws.Range("D80").Value = "Below Median"
